$d = $word.ActiveDocument

for ($i = 1; $i -le 10; $i++) {
    $old = "{{ANEXO_$i}}"
    $new = "{{%ANEXO_$i}}"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
